$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: convert existing is_active boolean-literal cells (G2:G67) into TRUE() formulas ---
for ($r = 2; $r -le 67; $r++) {
    $ws.Cells.Item($r, 7).Formula = "=TRUE()"
}

# --- Step 2: append 32 new rows (68-99) with Arabic (ara) translations ---
# row 68
$ws.Cells.Item(68, 1).Value = 10094
$ws.Cells.Item(68, 2).Value = 'bloodType'
$ws.Cells.Item(68, 3).Value = 'Blood Type'
$ws.Cells.Item(68, 4).Value = 'string'
$ws.Cells.Item(68, 5).Value = '
	{
		"code": "101",
		"value": "أ",
	}'
$ws.Cells.Item(68, 5).WrapText = $true
$ws.Cells.Item(68, 6).Value = 'ara'
$ws.Cells.Item(68, 7).Formula = "=TRUE()"
$ws.Cells.Item(68, 8).Value = 'superadmin'
$ws.Cells.Item(68, 9).Value = 'now()'
$ws.Rows.Item(68).RowHeight = 70.1

# row 69
$ws.Cells.Item(69, 1).Value = 10095
$ws.Cells.Item(69, 2).Value = 'bloodType'
$ws.Cells.Item(69, 3).Value = 'Blood Type'
$ws.Cells.Item(69, 4).Value = 'string'
$ws.Cells.Item(69, 5).Value = '
	{
		"code": "102",
		"value": "أ +",
	}'
$ws.Cells.Item(69, 5).WrapText = $true
$ws.Cells.Item(69, 6).Value = 'ara'
$ws.Cells.Item(69, 7).Formula = "=TRUE()"
$ws.Cells.Item(69, 8).Value = 'superadmin'
$ws.Cells.Item(69, 9).Value = 'now()'
$ws.Rows.Item(69).RowHeight = 70.1

# row 70
$ws.Cells.Item(70, 1).Value = 10096
$ws.Cells.Item(70, 2).Value = 'bloodType'
$ws.Cells.Item(70, 3).Value = 'Blood Type'
$ws.Cells.Item(70, 4).Value = 'string'
$ws.Cells.Item(70, 5).Value = '
	{
		"code": "103",
		"value": "أ-",
	}'
$ws.Cells.Item(70, 5).WrapText = $true
$ws.Cells.Item(70, 6).Value = 'ara'
$ws.Cells.Item(70, 7).Formula = "=TRUE()"
$ws.Cells.Item(70, 8).Value = 'superadmin'
$ws.Cells.Item(70, 9).Value = 'now()'
$ws.Rows.Item(70).RowHeight = 70.1

# row 71
$ws.Cells.Item(71, 1).Value = 10097
$ws.Cells.Item(71, 2).Value = 'bloodType'
$ws.Cells.Item(71, 3).Value = 'Blood Type'
$ws.Cells.Item(71, 4).Value = 'string'
$ws.Cells.Item(71, 5).Value = '
	{
		"code": "104",
		"value": "ب",
	}'
$ws.Cells.Item(71, 5).WrapText = $true
$ws.Cells.Item(71, 6).Value = 'ara'
$ws.Cells.Item(71, 7).Formula = "=TRUE()"
$ws.Cells.Item(71, 8).Value = 'superadmin'
$ws.Cells.Item(71, 9).Value = 'now()'
$ws.Rows.Item(71).RowHeight = 70.1

# row 72
$ws.Cells.Item(72, 1).Value = 10098
$ws.Cells.Item(72, 2).Value = 'bloodType'
$ws.Cells.Item(72, 3).Value = 'Blood Type'
$ws.Cells.Item(72, 4).Value = 'string'
$ws.Cells.Item(72, 5).Value = '
	{
		"code": "105",
		"value": "ب +",
	}'
$ws.Cells.Item(72, 5).WrapText = $true
$ws.Cells.Item(72, 6).Value = 'ara'
$ws.Cells.Item(72, 7).Formula = "=TRUE()"
$ws.Cells.Item(72, 8).Value = 'superadmin'
$ws.Cells.Item(72, 9).Value = 'now()'
$ws.Rows.Item(72).RowHeight = 70.1

# row 73
$ws.Cells.Item(73, 1).Value = 10099
$ws.Cells.Item(73, 2).Value = 'bloodType'
$ws.Cells.Item(73, 3).Value = 'Blood Type'
$ws.Cells.Item(73, 4).Value = 'string'
$ws.Cells.Item(73, 5).Value = '
	{
		"code": "106",
		"value": "ب-",
	}'
$ws.Cells.Item(73, 5).WrapText = $true
$ws.Cells.Item(73, 6).Value = 'ara'
$ws.Cells.Item(73, 7).Formula = "=TRUE()"
$ws.Cells.Item(73, 8).Value = 'superadmin'
$ws.Cells.Item(73, 9).Value = 'now()'
$ws.Rows.Item(73).RowHeight = 70.1

# row 74
$ws.Cells.Item(74, 1).Value = 10100
$ws.Cells.Item(74, 2).Value = 'bloodType'
$ws.Cells.Item(74, 3).Value = 'Blood Type'
$ws.Cells.Item(74, 4).Value = 'string'
$ws.Cells.Item(74, 5).Value = '
	{
		"code": "107",
		"value": "AB",
	}'
$ws.Cells.Item(74, 5).WrapText = $true
$ws.Cells.Item(74, 6).Value = 'ara'
$ws.Cells.Item(74, 7).Formula = "=TRUE()"
$ws.Cells.Item(74, 8).Value = 'superadmin'
$ws.Cells.Item(74, 9).Value = 'now()'
$ws.Rows.Item(74).RowHeight = 68.65

# row 75
$ws.Cells.Item(75, 1).Value = 10101
$ws.Cells.Item(75, 2).Value = 'bloodType'
$ws.Cells.Item(75, 3).Value = 'Blood Type'
$ws.Cells.Item(75, 4).Value = 'string'
$ws.Cells.Item(75, 5).Value = '
	{
		"code": "108",
		"value": "AB+",
	}'
$ws.Cells.Item(75, 5).WrapText = $true
$ws.Cells.Item(75, 6).Value = 'ara'
$ws.Cells.Item(75, 7).Formula = "=TRUE()"
$ws.Cells.Item(75, 8).Value = 'superadmin'
$ws.Cells.Item(75, 9).Value = 'now()'
$ws.Rows.Item(75).RowHeight = 68.65

# row 76
$ws.Cells.Item(76, 1).Value = 10102
$ws.Cells.Item(76, 2).Value = 'bloodType'
$ws.Cells.Item(76, 3).Value = 'Blood Type'
$ws.Cells.Item(76, 4).Value = 'string'
$ws.Cells.Item(76, 5).Value = '
	{
		"code": "109",
		"value": "AB-",
	}'
$ws.Cells.Item(76, 5).WrapText = $true
$ws.Cells.Item(76, 6).Value = 'ara'
$ws.Cells.Item(76, 7).Formula = "=TRUE()"
$ws.Cells.Item(76, 8).Value = 'superadmin'
$ws.Cells.Item(76, 9).Value = 'now()'
$ws.Rows.Item(76).RowHeight = 68.65

# row 77
$ws.Cells.Item(77, 1).Value = 10103
$ws.Cells.Item(77, 2).Value = 'bloodType'
$ws.Cells.Item(77, 3).Value = 'Blood Type'
$ws.Cells.Item(77, 4).Value = 'string'
$ws.Cells.Item(77, 5).Value = '
	{
		"code": "110",
		"value": "ا",
	}'
$ws.Cells.Item(77, 5).WrapText = $true
$ws.Cells.Item(77, 6).Value = 'ara'
$ws.Cells.Item(77, 7).Formula = "=TRUE()"
$ws.Cells.Item(77, 8).Value = 'superadmin'
$ws.Cells.Item(77, 9).Value = 'now()'
$ws.Rows.Item(77).RowHeight = 70.1

# row 78
$ws.Cells.Item(78, 1).Value = 10104
$ws.Cells.Item(78, 2).Value = 'bloodType'
$ws.Cells.Item(78, 3).Value = 'Blood Type'
$ws.Cells.Item(78, 4).Value = 'string'
$ws.Cells.Item(78, 5).Value = '
	{
		"code": "111",
		"value": "O+"
	}'
$ws.Cells.Item(78, 5).WrapText = $true
$ws.Cells.Item(78, 6).Value = 'ara'
$ws.Cells.Item(78, 7).Formula = "=TRUE()"
$ws.Cells.Item(78, 8).Value = 'superadmin'
$ws.Cells.Item(78, 9).Value = 'now()'
$ws.Rows.Item(78).RowHeight = 68.65

# row 79
$ws.Cells.Item(79, 1).Value = 10105
$ws.Cells.Item(79, 2).Value = 'bloodType'
$ws.Cells.Item(79, 3).Value = 'Blood Type'
$ws.Cells.Item(79, 4).Value = 'string'
$ws.Cells.Item(79, 5).Value = '
	{
		"code": "112",
		"value": "س-"
	}'
$ws.Cells.Item(79, 5).WrapText = $true
$ws.Cells.Item(79, 6).Value = 'ara'
$ws.Cells.Item(79, 7).Formula = "=TRUE()"
$ws.Cells.Item(79, 8).Value = 'superadmin'
$ws.Cells.Item(79, 9).Value = 'now()'
$ws.Rows.Item(79).RowHeight = 70.1

# row 80
$ws.Cells.Item(80, 1).Value = 10106
$ws.Cells.Item(80, 2).Value = 'bloodType'
$ws.Cells.Item(80, 3).Value = 'Blood Type'
$ws.Cells.Item(80, 4).Value = 'string'
$ws.Cells.Item(80, 5).Value = '
	{
		"code": "113",
		"value": "لا اعرف"
	}'
$ws.Cells.Item(80, 5).WrapText = $true
$ws.Cells.Item(80, 6).Value = 'ara'
$ws.Cells.Item(80, 7).Formula = "=TRUE()"
$ws.Cells.Item(80, 8).Value = 'superadmin'
$ws.Cells.Item(80, 9).Value = 'now()'
$ws.Rows.Item(80).RowHeight = 70.1

# row 81
$ws.Cells.Item(81, 1).Value = 10107
$ws.Cells.Item(81, 2).Value = 'bloodType'
$ws.Cells.Item(81, 3).Value = 'Blood Type'
$ws.Cells.Item(81, 4).Value = 'string'
$ws.Cells.Item(81, 5).Value = '
	{
		"code": "114",
		"value": "غير قابل للتطبيق"
	}'
$ws.Cells.Item(81, 5).WrapText = $true
$ws.Cells.Item(81, 6).Value = 'ara'
$ws.Cells.Item(81, 7).Formula = "=TRUE()"
$ws.Cells.Item(81, 8).Value = 'superadmin'
$ws.Cells.Item(81, 9).Value = 'now()'
$ws.Rows.Item(81).RowHeight = 70.1

# row 82
$ws.Cells.Item(82, 1).Value = 10108
$ws.Cells.Item(82, 2).Value = 'maritalStatus'
$ws.Cells.Item(82, 3).Value = 'Marital Status'
$ws.Cells.Item(82, 4).Value = 'string'
$ws.Cells.Item(82, 5).Value = '
	{
		"code": "101",
		"value": "غير مرتبطة
"
	}'
$ws.Cells.Item(82, 5).WrapText = $true
$ws.Cells.Item(82, 6).Value = 'ara'
$ws.Cells.Item(82, 7).Formula = "=TRUE()"
$ws.Cells.Item(82, 8).Value = 'superadmin'
$ws.Cells.Item(82, 9).Value = 'now()'
$ws.Rows.Item(82).RowHeight = 83.55

# row 83
$ws.Cells.Item(83, 1).Value = 10109
$ws.Cells.Item(83, 2).Value = 'maritalStatus'
$ws.Cells.Item(83, 3).Value = 'Marital Status'
$ws.Cells.Item(83, 4).Value = 'string'
$ws.Cells.Item(83, 5).Value = '
	{
		"code": "102",
		"value": "زوجت"
	}'
$ws.Cells.Item(83, 5).WrapText = $true
$ws.Cells.Item(83, 6).Value = 'ara'
$ws.Cells.Item(83, 7).Formula = "=TRUE()"
$ws.Cells.Item(83, 8).Value = 'superadmin'
$ws.Cells.Item(83, 9).Value = 'now()'
$ws.Rows.Item(83).RowHeight = 70.1

# row 84
$ws.Cells.Item(84, 1).Value = 10110
$ws.Cells.Item(84, 2).Value = 'maritalStatus'
$ws.Cells.Item(84, 3).Value = 'Marital Status'
$ws.Cells.Item(84, 4).Value = 'string'
$ws.Cells.Item(84, 5).Value = '
	{
		"code": "103",
		"value": "الأرامل"
	}'
$ws.Cells.Item(84, 5).WrapText = $true
$ws.Cells.Item(84, 6).Value = 'ara'
$ws.Cells.Item(84, 7).Formula = "=TRUE()"
$ws.Cells.Item(84, 8).Value = 'superadmin'
$ws.Cells.Item(84, 9).Value = 'now()'
$ws.Rows.Item(84).RowHeight = 70.1

# row 85
$ws.Cells.Item(85, 1).Value = 10111
$ws.Cells.Item(85, 2).Value = 'maritalStatus'
$ws.Cells.Item(85, 3).Value = 'Marital Status'
$ws.Cells.Item(85, 4).Value = 'string'
$ws.Cells.Item(85, 5).Value = '
	{
		"code": "104",
		"value": "مطلق"
	}'
$ws.Cells.Item(85, 5).WrapText = $true
$ws.Cells.Item(85, 6).Value = 'ara'
$ws.Cells.Item(85, 7).Formula = "=TRUE()"
$ws.Cells.Item(85, 8).Value = 'superadmin'
$ws.Cells.Item(85, 9).Value = 'now()'
$ws.Rows.Item(85).RowHeight = 70.1

# row 86
$ws.Cells.Item(86, 1).Value = 10112
$ws.Cells.Item(86, 2).Value = 'maritalStatus'
$ws.Cells.Item(86, 3).Value = 'Marital Status'
$ws.Cells.Item(86, 4).Value = 'string'
$ws.Cells.Item(86, 5).Value = '
	{
		"code": "105",
		"value": "مطلق من الناحية القانونية
"
	}'
$ws.Cells.Item(86, 5).WrapText = $true
$ws.Cells.Item(86, 6).Value = 'ara'
$ws.Cells.Item(86, 7).Formula = "=TRUE()"
$ws.Cells.Item(86, 8).Value = 'superadmin'
$ws.Cells.Item(86, 9).Value = 'now()'
$ws.Rows.Item(86).RowHeight = 83.55

# row 87
$ws.Cells.Item(87, 1).Value = 10113
$ws.Cells.Item(87, 2).Value = 'maritalStatus'
$ws.Cells.Item(87, 3).Value = 'Marital Status'
$ws.Cells.Item(87, 4).Value = 'string'
$ws.Cells.Item(87, 5).Value = '
	{
		"code": "106",
		"value": "ملغاة"
	}'
$ws.Cells.Item(87, 5).WrapText = $true
$ws.Cells.Item(87, 6).Value = 'ara'
$ws.Cells.Item(87, 7).Formula = "=TRUE()"
$ws.Cells.Item(87, 8).Value = 'superadmin'
$ws.Cells.Item(87, 9).Value = 'now()'
$ws.Rows.Item(87).RowHeight = 70.1

# row 88
$ws.Cells.Item(88, 1).Value = 10114
$ws.Cells.Item(88, 2).Value = 'maritalStatus'
$ws.Cells.Item(88, 3).Value = 'Marital Status'
$ws.Cells.Item(88, 4).Value = 'string'
$ws.Cells.Item(88, 5).Value = '
	{
		"code": "107",
		"value": "مبطل"
	}'
$ws.Cells.Item(88, 5).WrapText = $true
$ws.Cells.Item(88, 6).Value = 'ara'
$ws.Cells.Item(88, 7).Formula = "=TRUE()"
$ws.Cells.Item(88, 8).Value = 'superadmin'
$ws.Cells.Item(88, 9).Value = 'now()'
$ws.Rows.Item(88).RowHeight = 70.1

# row 89
$ws.Cells.Item(89, 1).Value = 10115
$ws.Cells.Item(89, 2).Value = 'registrationType'
$ws.Cells.Item(89, 3).Value = 'Registration Type'
$ws.Cells.Item(89, 4).Value = 'string'
$ws.Cells.Item(89, 5).Value = '
	{
		"code": "Document-based",
		"value": "مستند إلى المستند"
	}'
$ws.Cells.Item(89, 5).WrapText = $true
$ws.Cells.Item(89, 6).Value = 'ara'
$ws.Cells.Item(89, 7).Formula = "=TRUE()"
$ws.Cells.Item(89, 8).Value = 'superadmin'
$ws.Cells.Item(89, 9).Value = 'now()'
$ws.Rows.Item(89).RowHeight = 70.1

# row 90
$ws.Cells.Item(90, 1).Value = 10116
$ws.Cells.Item(90, 2).Value = 'registrationType'
$ws.Cells.Item(90, 3).Value = 'Registration Type'
$ws.Cells.Item(90, 4).Value = 'string'
$ws.Cells.Item(90, 5).Value = '
	{
		"code": "Introducer-based",
		"value": "المعرف القائم"
	}'
$ws.Cells.Item(90, 5).WrapText = $true
$ws.Cells.Item(90, 6).Value = 'ara'
$ws.Cells.Item(90, 7).Formula = "=TRUE()"
$ws.Cells.Item(90, 8).Value = 'superadmin'
$ws.Cells.Item(90, 9).Value = 'now()'
$ws.Rows.Item(90).RowHeight = 70.1

# row 91
$ws.Cells.Item(91, 1).Value = 10117
$ws.Cells.Item(91, 2).Value = 'modeOfClaim'
$ws.Cells.Item(91, 3).Value = 'Mode of Claim'
$ws.Cells.Item(91, 4).Value = 'string'
$ws.Cells.Item(91, 5).Value = '
	{
		"code": "101",
		"value": "يلتقط"
	}'
$ws.Cells.Item(91, 5).WrapText = $true
$ws.Cells.Item(91, 6).Value = 'ara'
$ws.Cells.Item(91, 7).Formula = "=TRUE()"
$ws.Cells.Item(91, 8).Value = 'superadmin'
$ws.Cells.Item(91, 9).Value = 'now()'
$ws.Rows.Item(91).RowHeight = 70.1

# row 92
$ws.Cells.Item(92, 1).Value = 10118
$ws.Cells.Item(92, 2).Value = 'modeOfClaim'
$ws.Cells.Item(92, 3).Value = 'Mode of Claim'
$ws.Cells.Item(92, 4).Value = 'string'
$ws.Cells.Item(92, 5).Value = '
	{
		"code": "102",
		"value": "التسليم إلى العنوان الدائم"
	}'
$ws.Cells.Item(92, 5).WrapText = $true
$ws.Cells.Item(92, 6).Value = 'ara'
$ws.Cells.Item(92, 7).Formula = "=TRUE()"
$ws.Cells.Item(92, 8).Value = 'superadmin'
$ws.Cells.Item(92, 9).Value = 'now()'
$ws.Rows.Item(92).RowHeight = 70.1

# row 93
$ws.Cells.Item(93, 1).Value = 10119
$ws.Cells.Item(93, 2).Value = 'modeOfClaim'
$ws.Cells.Item(93, 3).Value = 'Mode of Claim'
$ws.Cells.Item(93, 4).Value = 'string'
$ws.Cells.Item(93, 5).Value = '
	{
		"code": "103",
		"value": "التسليم إلى العنوان الحالي"
	}'
$ws.Cells.Item(93, 5).WrapText = $true
$ws.Cells.Item(93, 6).Value = 'ara'
$ws.Cells.Item(93, 7).Formula = "=TRUE()"
$ws.Cells.Item(93, 8).Value = 'superadmin'
$ws.Cells.Item(93, 9).Value = 'now()'
$ws.Rows.Item(93).RowHeight = 70.1

# row 94
$ws.Cells.Item(94, 1).Value = 10120
$ws.Cells.Item(94, 2).Value = 'gender'
$ws.Cells.Item(94, 3).Value = 'Gender'
$ws.Cells.Item(94, 4).Value = 'string'
$ws.Cells.Item(94, 5).Value = '
	{
		"code": "MLE",
		"value": "ذكر"
	}'
$ws.Cells.Item(94, 5).WrapText = $true
$ws.Cells.Item(94, 6).Value = 'ara'
$ws.Cells.Item(94, 7).Formula = "=TRUE()"
$ws.Cells.Item(94, 8).Value = 'superadmin'
$ws.Cells.Item(94, 9).Value = 'now()'
$ws.Rows.Item(94).RowHeight = 70.1

# row 95
$ws.Cells.Item(95, 1).Value = 10121
$ws.Cells.Item(95, 2).Value = 'gender'
$ws.Cells.Item(95, 3).Value = 'Gender'
$ws.Cells.Item(95, 4).Value = 'string'
$ws.Cells.Item(95, 5).Value = '
	{
		"code": "FLE",
		"value": "أنثى"
	}'
$ws.Cells.Item(95, 5).WrapText = $true
$ws.Cells.Item(95, 6).Value = 'ara'
$ws.Cells.Item(95, 7).Formula = "=TRUE()"
$ws.Cells.Item(95, 8).Value = 'superadmin'
$ws.Cells.Item(95, 9).Value = 'now()'
$ws.Rows.Item(95).RowHeight = 70.1

# row 96
$ws.Cells.Item(96, 1).Value = 10122
$ws.Cells.Item(96, 2).Value = 'gender'
$ws.Cells.Item(96, 3).Value = 'Gender'
$ws.Cells.Item(96, 4).Value = 'string'
$ws.Cells.Item(96, 5).Value = '
	{
		"code": "OTH",
		"value": "آحرون"
	}'
$ws.Cells.Item(96, 5).WrapText = $true
$ws.Cells.Item(96, 6).Value = 'ara'
$ws.Cells.Item(96, 7).Formula = "=TRUE()"
$ws.Cells.Item(96, 8).Value = 'superadmin'
$ws.Cells.Item(96, 9).Value = 'now()'
$ws.Rows.Item(96).RowHeight = 70.1

# row 97
$ws.Cells.Item(97, 1).Value = 10123
$ws.Cells.Item(97, 2).Value = 'residenceStatus'
$ws.Cells.Item(97, 3).Value = 'residenceStatus'
$ws.Cells.Item(97, 4).Value = 'string'
$ws.Cells.Item(97, 5).Value = '
	{
		"code": "FR",
		"value": "أجنبي"
	}'
$ws.Cells.Item(97, 5).WrapText = $true
$ws.Cells.Item(97, 6).Value = 'ara'
$ws.Cells.Item(97, 7).Formula = "=TRUE()"
$ws.Cells.Item(97, 8).Value = 'superadmin'
$ws.Cells.Item(97, 9).Value = 'now()'
$ws.Rows.Item(97).RowHeight = 70.1

# row 98
$ws.Cells.Item(98, 1).Value = 10124
$ws.Cells.Item(98, 2).Value = 'residenceStatus'
$ws.Cells.Item(98, 3).Value = 'residenceStatus'
$ws.Cells.Item(98, 4).Value = 'string'
$ws.Cells.Item(98, 5).Value = '
	{
		"code": "NFR",
		"value": "غير أجنبي"
	}'
$ws.Cells.Item(98, 5).WrapText = $true
$ws.Cells.Item(98, 6).Value = 'ara'
$ws.Cells.Item(98, 7).Formula = "=TRUE()"
$ws.Cells.Item(98, 8).Value = 'superadmin'
$ws.Cells.Item(98, 9).Value = 'now()'
$ws.Rows.Item(98).RowHeight = 70.1

# row 99
$ws.Cells.Item(99, 1).Value = 10125
$ws.Cells.Item(99, 2).Value = 'preferredLang'
$ws.Cells.Item(99, 3).Value = 'user preferred Language'
$ws.Cells.Item(99, 4).Value = 'string'
$ws.Cells.Item(99, 5).Value = '{"value":"عربى","code":"ara"}'
$ws.Cells.Item(99, 5).WrapText = $true
$ws.Cells.Item(99, 6).Value = 'ara'
$ws.Cells.Item(99, 7).Formula = "=TRUE()"
$ws.Cells.Item(99, 8).Value = 'superadmin'
$ws.Cells.Item(99, 9).Value = 'now()'
$ws.Rows.Item(99).RowHeight = 28.35

# row 99 additionally uses wrap-text style on B,C,D,F,H,I (mirrors source xf s="2")
$ws.Cells.Item(99, 2).WrapText = $true
$ws.Cells.Item(99, 3).WrapText = $true
$ws.Cells.Item(99, 4).WrapText = $true
$ws.Cells.Item(99, 6).WrapText = $true
$ws.Cells.Item(99, 8).WrapText = $true
$ws.Cells.Item(99, 9).WrapText = $true

# --- Step 3: misc view/formatting tweaks ---
$ws.Columns.Item(1).ColumnWidth = 11.7
$ws.Range("E79").Select() | Out-Null
